$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: new columns AD, AE, AF (Wins, Losses, Ties)
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy the header style from AC1 (which already has the bold/border/center style) to AD1:AF1
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122) # xlPasteFormats

for ($r = 2; $r -le 53; $r++) {
    $ws.Cells.Item($r, 30).Value = 75   # AD
    $ws.Cells.Item($r, 31).Value = 87   # AE
    $ws.Cells.Item($r, 32).Value = 0    # AF
}
